# Applies the "cryptos list" refresh described in the commit:
#   "Updated cryptos list on Wed Jan 31 10:33:49 UTC 2024 with GitHub Actions"
#
# Every row is rewritten with refreshed Price/Volume(1h) figures; rows 21-51 also
# shift up by one slot (a new coin, RocketPoolETH, enters at the bottom) while the
# numeric rank in column A stays put. All target cells are plain text (inline/shared
# strings in the original file), so numeric-looking prices are written with the cell
# NumberFormat forced to Text ('@') first -- otherwise Excel would silently coerce
# values like "0.0790" or "6.10" into numbers and mangle the significant trailing
# zeros / drop the string formatting entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r = A1-style cell reference, v = new text value
$updates = @(
    @{ r = 'D2'; v = '42.748.05' }
    @{ r = 'E2'; v = '  -1.47%  ' }
    @{ r = 'D3'; v = '2.306.40' }
    @{ r = 'E3'; v = '  -0.17%  ' }
    @{ r = 'E4'; v = '  +0.01%  ' }
    @{ r = 'D5'; v = '304.22' }
    @{ r = 'E5'; v = '  -2.13%  ' }
    @{ r = 'D6'; v = '99.43' }
    @{ r = 'E6'; v = '  -4.22%  ' }
    @{ r = 'E7'; v = '  -5.09%  ' }
    @{ r = 'E9'; v = '  -5.28%  ' }
    @{ r = 'E10'; v = '  -6.19%  ' }
    @{ r = 'D11'; v = '51.86' }
    @{ r = 'E11'; v = '  -1.72%  ' }
    @{ r = 'D12'; v = '0.0790' }
    @{ r = 'E12'; v = '  -2.78%  ' }
    @{ r = 'E13'; v = '  +0.64%  ' }
    @{ r = 'E14'; v = '  -3.45%  ' }
    @{ r = 'D15'; v = '2.664.58' }
    @{ r = 'E15'; v = '  -0.04%  ' }
    @{ r = 'D16'; v = '15.72' }
    @{ r = 'E16'; v = '  +4.06%  ' }
    @{ r = 'D17'; v = '2.305.56' }
    @{ r = 'E17'; v = '  -0.19%  ' }
    @{ r = 'D18'; v = '0.822' }
    @{ r = 'E18'; v = '  +1.50%  ' }
    @{ r = 'D19'; v = '42.657.20' }
    @{ r = 'D20'; v = '0.0₃0902' }
    @{ r = 'E20'; v = '  -2.61%  ' }
    @{ r = 'B21'; v = 'InternetComputer(DFINITY)' }
    @{ r = 'C21'; v = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ r = 'D21'; v = '11.56' }
    @{ r = 'E21'; v = '  -5.08%  ' }
    @{ r = 'B22'; v = 'Uniswap' }
    @{ r = 'C22'; v = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ r = 'D22'; v = '6.10' }
    @{ r = 'E22'; v = '  -1.24%  ' }
    @{ r = 'D23'; v = '69.16' }
    @{ r = 'E23'; v = '  +1.57%  ' }
    @{ r = 'D24'; v = '234.85' }
    @{ r = 'E24'; v = '  -3.16%  ' }
    @{ r = 'E25'; v = '  -2.13%  ' }
    @{ r = 'E26'; v = '  -3.31%  ' }
    @{ r = 'E27'; v = '  -0.08%  ' }
    @{ r = 'B28'; v = 'EthereumClassic' }
    @{ r = 'C28'; v = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ r = 'D28'; v = '25.26' }
    @{ r = 'E28'; v = '  +1.29%  ' }
    @{ r = 'B29'; v = 'Toncoin' }
    @{ r = 'C29'; v = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ r = 'D29'; v = '2.17' }
    @{ r = 'E29'; v = '  -5.51%  ' }
    @{ r = 'D30'; v = '34.66' }
    @{ r = 'E30'; v = '  -6.46%  ' }
    @{ r = 'B31'; v = 'Cosmos' }
    @{ r = 'C31'; v = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ r = 'D31'; v = '9.20' }
    @{ r = 'E31'; v = '  -4.59%  ' }
    @{ r = 'B32'; v = 'Monero' }
    @{ r = 'C32'; v = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ r = 'D32'; v = '161.95' }
    @{ r = 'E32'; v = '  -3.39%  ' }
    @{ r = 'B33'; v = 'FirstDigitalUSD' }
    @{ r = 'C33'; v = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' }
    @{ r = 'D33'; v = '0.999' }
    @{ r = 'E33'; v = '  -0.02%  ' }
    @{ r = 'B34'; v = 'Filecoin' }
    @{ r = 'C34'; v = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ r = 'D34'; v = '5.05' }
    @{ r = 'E34'; v = '  -4.24%  ' }
    @{ r = 'B35'; v = 'RenderToken' }
    @{ r = 'C35'; v = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ r = 'D35'; v = '4.63' }
    @{ r = 'E35'; v = '  +4.09%  ' }
    @{ r = 'B36'; v = 'WEMIXToken' }
    @{ r = 'C36'; v = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ r = 'D36'; v = '2.44' }
    @{ r = 'E36'; v = '  -3.53%  ' }
    @{ r = 'B37'; v = 'Hedera' }
    @{ r = 'C37'; v = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ r = 'D37'; v = '0.0717' }
    @{ r = 'E37'; v = '  -3.63%  ' }
    @{ r = 'B38'; v = 'Celestia' }
    @{ r = 'C38'; v = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia' }
    @{ r = 'D38'; v = '17.03' }
    @{ r = 'E38'; v = '  -7.24%  ' }
    @{ r = 'B39'; v = 'LidoDAOToken' }
    @{ r = 'C39'; v = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ r = 'D39'; v = '2.89' }
    @{ r = 'E39'; v = '  -5.28%  ' }
    @{ r = 'B40'; v = 'ARBITRUM' }
    @{ r = 'C40'; v = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ r = 'D40'; v = '1.81' }
    @{ r = 'E40'; v = '  -3.84%  ' }
    @{ r = 'B41'; v = 'Kaspa' }
    @{ r = 'C41'; v = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ r = 'D41'; v = '0.101' }
    @{ r = 'E41'; v = '  -5.05%  ' }
    @{ r = 'B42'; v = 'Stellar' }
    @{ r = 'C42'; v = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ r = 'D42'; v = '0.112' }
    @{ r = 'E42'; v = '  -3.43%  ' }
    @{ r = 'B43'; v = 'ApeXProtocol' }
    @{ r = 'C43'; v = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
    @{ r = 'D43'; v = '2.46' }
    @{ r = 'E43'; v = '  -8.99%  ' }
    @{ r = 'B44'; v = 'Maker' }
    @{ r = 'C44'; v = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ r = 'D44'; v = '1.993.10' }
    @{ r = 'E44'; v = '  +0.23%  ' }
    @{ r = 'B45'; v = 'EnergySwap' }
    @{ r = 'C45'; v = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ r = 'D45'; v = '18.74' }
    @{ r = 'E45'; v = '  -1.63%  ' }
    @{ r = 'B46'; v = 'VeChain' }
    @{ r = 'C46'; v = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ r = 'D46'; v = '0.0280' }
    @{ r = 'E46'; v = '  -4.29%  ' }
    @{ r = 'B47'; v = 'FraxShare' }
    @{ r = 'C47'; v = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ r = 'D47'; v = '10.26' }
    @{ r = 'E47'; v = '  +2.54%  ' }
    @{ r = 'B48'; v = 'NEARProtocol' }
    @{ r = 'C48'; v = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ r = 'D48'; v = '2.87' }
    @{ r = 'E48'; v = '  -5.99%  ' }
    @{ r = 'B49'; v = 'MultiversX' }
    @{ r = 'C49'; v = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' }
    @{ r = 'D49'; v = '55.46' }
    @{ r = 'E49'; v = '  -0.54%  ' }
    @{ r = 'B50'; v = 'HuobiToken' }
    @{ r = 'C50'; v = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ r = 'D50'; v = '2.87' }
    @{ r = 'E50'; v = '  -2.79%  ' }
    @{ r = 'B51'; v = 'RocketPoolETH' }
    @{ r = 'C51'; v = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' }
    @{ r = 'D51'; v = '2.534.49' }
    @{ r = 'E51'; v = '  -0.02%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.r)
    if ($u.v -match "^[+-]?\d+(\.\d+)?$") {
        # Looks like a number (e.g. "51.86"/"0.0790") but must stay text, matching
        # the source inline-string cell -- force Text format before assigning.
        $cell.NumberFormat = '@'
    }
    $cell.Value = $u.v
}
